$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM row 28: TVS / Schottky diodes added on Dyna data lines (D4-D7)
$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "SMAJ5920"
$ws.Range("C28").Value = "DIODE-SCHOTTKY-MBRA140"
$ws.Range("D28").Value = "SMA-DIODE"
$ws.Range("E28").Value = "D4, D5, D6, D7"
$ws.Range("F28").Value = "Schottky diodes in SFE's production catalog"
$ws.Range("G28").Value = "Y"
$ws.Range("J28").Value = "51T4985"

# Match row height used throughout the rest of the sheet
$ws.Rows.Item(28).RowHeight = 15.75

# Qty..Description get the bordered "SUPPLIER"-style look (Calibri 12, black, boxed)
$dataRange = $ws.Range("A28:F28")
$dataRange.Font.Name = "Calibri"
$dataRange.Font.Size = 12
$dataRange.Font.Color = 0
$dataRange.Borders.LineStyle = 1

# "Needs Ordered" flag keeps the normal font but vertically centered
$ws.Range("G28").Font.Name = "Arial"
$ws.Range("G28").Font.Size = 10
$ws.Range("G28").VerticalAlignment = -4108

# Move the active selection like the author left it
$ws.Range("G21").Select()
